$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "이항 분포"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/04/23/binomial_distribution.html"

$ws.Range("D9").Value = "[공지] SBS-Pabii 대학원 관련 공지"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/notice-kakaotalk-chat-service/#utm_source=rss&utm_medium=rss&utm_campaign=notice-kakaotalk-chat-service"

$ws.Range("D16").Value = "Interpretable and fine-grained visual explanations for CNNs 내용 정리 [XAI-7]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/149"

$ws.Range("D42").Value = "닷넷 3.5 오프라인 설치 방법"
$ws.Range("E42").Value = "https://kjk92.tistory.com/68"

$ws.Range("D51").Value = "MySQL workbench에서 delete 안 될 때 해결방법  (error code: 1175)"
$ws.Range("E51").Value = "https://bskyvision.com/1181"
